# Script: apply "extraction conversion scores" edit to Table 1 worksheet.
# - Splits the combined "PanNETs"/"PanNECs" breakdown (in columns I/J of row 15-16)
#   into a proper additional table row (row 16), matching the other rows' layout.
# - Row 15 (PanNETs) becomes the "Pancreatic Neuroendocrine Tumor" row, now also
#   carrying a Tissue type ("unknown") and its count moves into column H.
# - Row 16 (PanNECs) becomes a new explicit row for "Pancreatic Neuroendocrine
#   Carcinoma" with the same Dataset/Chip design/Tissue as row 15.
# - The now unused columns I and J are cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15: update Tumor type, add Tissue type, move count into H ---
$ws.Range("E15").Value = "Pancreatic Neuroendocrine Tumor"
$ws.Range("G15").Value = "unknown"
$ws.Range("H15").Value = 30

# --- Row 16: new explicit data row for Pancreatic Neuroendocrine Carcinoma ---
$ws.Range("C16").Value = "JGAS000359"
$ws.Range("D16").Value = "EPIC"
$ws.Range("E16").Value = "Pancreatic Neuroendocrine Carcinoma"
$ws.Range("F16").Value = "Primary"
$ws.Range("G16").Value = "unknown"
$ws.Range("H16").Value = 14

# --- Clear now-unused columns I and J on rows 15 and 16 ---
$ws.Range("I15:J16").Clear()

# --- Match row 16's E/F formatting to row 15's (same style class as other data rows) ---
$ws.Range("E15:F15").Copy()
$ws.Range("E16:F16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Selection cosmetics (matches author's saved selection) ---
$ws.Range("M21").Select()
